# Apply cryptocurrency price/volume updates (GitHub Actions scheduled refresh)
# Note: Price values in column D are stored as literal text (the source data
# uses "." as a thousands-style separator, e.g. "29.425.92", so it can't be a
# real number). For D values that look like a plain decimal number (single
# "."), a leading apostrophe is used to force Excel to keep them as text
# instead of auto-converting to a Number (which would drop trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.425.92'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '1.850.32'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''240.34'
$ws.Range('D6').Value = '''0.6294'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '''0.07647'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = '''24.93'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('D11').Value = '2.077.86'
$ws.Range('E11').Value = '  +12.30%  '
$ws.Range('D12').Value = '''0.07743'
$ws.Range('D13').Value = '''5.035'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').Value = '''0.6812'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').Value = '''0.00001063'
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').Value = '''83.43'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '''6.182'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '29.518.21'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').Value = '''228.86'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('D22').Value = '''7.463'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '''157.59'
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('D26').Value = '''8.438'
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').Value = '''1.386'
$ws.Range('E28').Value = '  +6.66%  '
$ws.Range('D29').Value = '''1.466'
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').Value = '''0.05615'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').Value = '''4.132'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').Value = '''4.050'
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('D33').Value = '''1.847'
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('D35').Value = '''0.7005'
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = '''0.01804'
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('D38').Value = '1.229.46'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('D39').Value = '''2.716'
$ws.Range('E39').Value = '  -2.15%  '
$ws.Range('D40').Value = '''6.445'
$ws.Range('E40').Value = '  +0.55%  '
$ws.Range('D41').Value = '''0.9084'
$ws.Range('E41').Value = '  +0.30%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').Value = '''101.90'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').Value = '''66.10'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').Value = '''0.00000000120'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('D46').Value = '''7.194'
$ws.Range('E46').Value = '  +0.28%  '
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').Value = '''0.1155'
$ws.Range('E48').Value = '  +3.27%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''1.686'
$ws.Range('E49').Value = '  +0.59%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''8.993'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').Value = '''0.05708'
$ws.Range('E51').Value = '  -0.01%  '
